$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Column P: "Meets Both - AND()" formulas for rows 3..12 ---
$ws.Range("P3:P12").Formula = "=AND(L3,M3)"

# A stray space value typed (and later left) in Q11
$ws.Range("Q11").Value = " "

# --- Summary rows 14-17 ---

# Row 14: Average (AVERAGE())
$ws.Range("C14").Formula = "=AVERAGE(C3:C12)"
$ws.Range("D14:E14").Formula = "=AVERAGE(D3:D12)"
$ws.Range("D14").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# Row 15: Total Count
$ws.Range("C15").Formula = "=COUNT(C3:C12)"
$ws.Range("D15:E15").Formula = "=COUNT(D3:D12)"
$ws.Range("D15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# Stray formatted (but empty) cell left over from navigating around
$ws.Range("R15").NumberFormat = "mmm-yy"

# Row 16: Meets Goals Count (hard-coded thresholds)
$ws.Range("C16").Formula = '=COUNTIF($C$3:$C$12,"<=5")'
$ws.Range("D16").Formula = '=COUNTIF($D$3:$D$12,">=90000")'
$ws.Range("D16").Copy()
$ws.Range("E16").PasteSpecial(-4122)

# Row 17: Meets Goals Count, copied format down from row 16 first
$ws.Range("C16:D16").Copy()
$ws.Range("C17:D17").PasteSpecial(-4122)
$ws.Range("C17").Formula = '=COUNTIF($C$3:$C$12,"<="&$S$3)'
$ws.Range("D17").Formula = '=COUNTIF($D$3:$D$12,">=90000")'

$ws.Range("D17").Select() | Out-Null
